# Update column F (dSF) values for the rows that were re-pulled / recalculated.
# Mapping: spreadsheet row number -> new value (per the diff of the source data).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = 1
    5  = -1
    9  = -1
    10 = 1
    12 = -5
    16 = -3
    17 = 4
    22 = 6
    23 = 0
    24 = 1
    26 = 4
    29 = -3
    30 = 1
    36 = 2
    38 = 3
    39 = 0
    41 = 2
    43 = -7
    53 = -2
    57 = 3
    58 = -4
    62 = -3
    64 = -3
    65 = 3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
